$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 713408145429.4626
    3  = 242814620310.7914
    4  = 35356419041.60538
    5  = 31952860863.86129
    6  = 23219394016.70679
    7  = 13054128698.86161
    8  = 10572957351.13158
    9  = 8855656019.970976
    10 = 8207533770.940097
    11 = 7899877325.028001
    12 = 7172856209.641662
    13 = 7052109471.634726
    14 = 6549647879.810579
    15 = 5943020494.124944
    16 = 5001855713.450159
    17 = 4696164014.753898
    18 = 4300199471.765975
    19 = 3734294537.424081
    20 = 3461093890.421829
    21 = 3254968381.180032
    22 = 3233171860.471212
    23 = 2978080122.785878
    24 = 2855259008.35902
    25 = 2685551641.418438
    26 = 2406688509.304906
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
